$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 166. The rows currently occupying
# 166-177 shift down to 168-179 (carrying their values/formatting).
$ws.Range("A166:A167").EntireRow.Insert()

# New row 166: Maracuya "Especial" quality, week of 2023-05-29
$ws.Cells.Item(166, 1).Value = 1
$ws.Cells.Item(166, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(166, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(166, 4).Value = "2023-05-29"
$ws.Cells.Item(166, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(166, 5).Value = 15
$ws.Cells.Item(166, 6).Value = "Fruta"
$ws.Cells.Item(166, 7).Value = 100108
$ws.Cells.Item(166, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(166, 9).Value = 100108003
$ws.Cells.Item(166, 10).Value = "Maracuyá"
$ws.Cells.Item(166, 11).Value = "Sin especificar"
$ws.Cells.Item(166, 12).Value = "Especial"
$ws.Cells.Item(166, 13).Value = 100
$ws.Cells.Item(166, 14).Value = 34000
$ws.Cells.Item(166, 15).Value = 35000
$ws.Cells.Item(166, 16).Value = 34500
$ws.Cells.Item(166, 17).Value = "`$/caja 20 kilos"
$ws.Cells.Item(166, 18).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(166, 19).Value = 1725
$ws.Cells.Item(166, 20).Value = 20

# New row 167: Maracuya "Primera" quality, week of 2023-05-29
$ws.Cells.Item(167, 1).Value = 1
$ws.Cells.Item(167, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(167, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(167, 4).Value = "2023-05-29"
$ws.Cells.Item(167, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(167, 5).Value = 15
$ws.Cells.Item(167, 6).Value = "Fruta"
$ws.Cells.Item(167, 7).Value = 100108
$ws.Cells.Item(167, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(167, 9).Value = 100108003
$ws.Cells.Item(167, 10).Value = "Maracuyá"
$ws.Cells.Item(167, 11).Value = "Sin especificar"
$ws.Cells.Item(167, 12).Value = "Primera"
$ws.Cells.Item(167, 13).Value = 130
$ws.Cells.Item(167, 14).Value = 30000
$ws.Cells.Item(167, 15).Value = 31000
$ws.Cells.Item(167, 16).Value = 30500
$ws.Cells.Item(167, 17).Value = "`$/caja 20 kilos"
$ws.Cells.Item(167, 18).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(167, 19).Value = 1525
$ws.Cells.Item(167, 20).Value = 20
